$wb = $excel.ActiveWorkbook
Write-Output "Sheets:"
foreach ($ws in $wb.Worksheets) {
    Write-Output $ws.Name
}
